# Applies the 16:18:55 re-scrape update for Línea 141 (commit: "Horarios actualizados Línea 141 - 281").
# For each sheet: refresh the "Última actualización" / "Total filas" header cells, then rewrite the full
# data body (rows 6+) with the freshly-scraped + re-sorted (by Hora_Llegada) row set.
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 16:18:55"
$ws.Range("A3").Value = "Total filas: 353"
$ws.Range("A6:E358").ClearContents()
$rows = @(
    @("04:18:53", "04:45", "215A_EL PATO", 27, "LP1912"),
    @("04:40:33", "04:46", "215A_EL PATO", 6, "LP1912"),
    @("04:18:53", "04:53", "11_ETCHEVERRY", 35, "LP1912"),
    @("04:54:03", "04:55", "11_ETCHEVERRY", 1, "LP1912"),
    @("04:54:03", "05:15", "14_ABASTO", 21, "LP1912"),
    @("04:18:53", "05:16", "17_ROMERO", 58, "LP1912"),
    @("04:40:33", "05:16", "14_ABASTO", 36, "LP1912"),
    @("05:18:56", "05:20", "14_ABASTO", 2, "LP1912"),
    @("04:18:53", "05:21", "23_HERNANDEZ", 63, "LP1912"),
    @("04:40:33", "05:22", "23_HERNANDEZ", 42, "LP1912"),
    @("04:18:53", "05:34", "215B_EL PATO", 76, "LP1912"),
    @("04:54:03", "05:35", "215B_EL PATO", 41, "LP1912"),
    @("04:18:53", "05:37", "14_ABASTO", 79, "LP1912"),
    @("04:18:53", "05:46", "15_ABASTO", 88, "LP1912"),
    @("05:49:40", "05:51", "14_ABASTO", 2, "LP1912"),
    @("04:40:33", "06:04", "16_SANTA ANA", 84, "LP1912"),
    @("05:49:40", "06:06", "16_SANTA ANA", 17, "LP1912"),
    @("04:18:53", "06:07", "16_SANTA ANA", 109, "LP1912"),
    @("04:18:53", "06:11", "215A_EL PATO", 113, "LP1912"),
    @("04:54:03", "06:12", "215A_EL PATO", 78, "LP1912"),
    @("04:18:53", "06:13", "225_HARAS DEL SUR", 115, "LP1912"),
    @("04:40:33", "06:14", "225_HARAS DEL SUR", 94, "LP1912"),
    @("06:15:04", "06:15", "225_HARAS DEL SUR", 0, "LP1912"),
    @("05:18:56", "06:20", "26_HERNANDEZ", 62, "LP1912"),
    @("04:40:33", "06:21", "26_HERNANDEZ", 101, "LP1912"),
    @("05:18:56", "06:26", "23_HERNANDEZ", 68, "LP1912"),
    @("04:40:33", "06:27", "23_HERNANDEZ", 107, "LP1912"),
    @("06:15:04", "06:28", "23_HERNANDEZ", 13, "LP1912"),
    @("04:40:33", "06:29", "86_EST CHICA-ESC AGRARIA", 109, "LP1912"),
    @("04:54:03", "06:30", "86_EST CHICA-ESC AGRARIA", 96, "LP1912"),
    @("04:40:33", "06:31", "16_SANTA ANA", 111, "LP1912"),
    @("05:18:56", "06:43", "225_C ROCA-H SUR", 85, "LP1912"),
    @("04:54:03", "06:44", "225_C ROCA-H SUR", 110, "LP1912"),
    @("05:18:56", "06:46", "215C_EL PATO", 88, "LP1912"),
    @("06:43:40", "06:46", "225_C ROCA-H SUR", 3, "LP1912"),
    @("04:54:03", "06:47", "215C_EL PATO", 113, "LP1912"),
    @("05:18:56", "06:58", "10_OLMOS", 100, "LP1912"),
    @("06:57:30", "06:58", "14_ABASTO", 1, "LP1912"),
    @("05:18:56", "06:59", "14_ABASTO", 101, "LP1912"),
    @("06:15:04", "07:00", "14_ABASTO", 45, "LP1912"),
    @("06:15:04", "07:01", "16_SANTA ANA", 46, "LP1912"),
    @("05:49:40", "07:04", "23_HERNANDEZ", 75, "LP1912"),
    @("05:18:56", "07:04", "15_ABASTO", 106, "LP1912"),
    @("05:49:40", "07:05", "15_ABASTO", 76, "LP1912"),
    @("05:18:56", "07:06", "225_GOMEZ", 108, "LP1912"),
    @("05:49:40", "07:07", "225_GOMEZ", 78, "LP1912"),
    @("05:18:56", "07:11", "215A_EL PATO", 113, "LP1912"),
    @("06:15:04", "07:12", "215A_EL PATO", 57, "LP1912"),
    @("05:18:56", "07:15", "11_ETCHEVERRY", 117, "LP1912"),
    @("06:43:40", "07:16", "16_SANTA ANA", 33, "LP1912"),
    @("06:15:04", "07:16", "11_ETCHEVERRY", 61, "LP1912"),
    @("07:20:40", "07:20", "10_OLMOS", 0, "LP1912"),
    @("05:49:40", "07:21", "26_HERNANDEZ", 92, "LP1912"),
    @("06:15:04", "07:23", "10_OLMOS", 68, "LP1912"),
    @("05:49:40", "07:29", "10_OLMOS", 100, "LP1912"),
    @("05:49:40", "07:31", "11_ETCHEVERRY", 102, "LP1912"),
    @("06:15:04", "07:32", "11_ETCHEVERRY", 77, "LP1912"),
    @("05:49:40", "07:32", "84_COLONIA URQUIZA-ESC 49", 103, "LP1912"),
    @("07:20:40", "07:34", "23_HERNANDEZ", 14, "LP1912"),
    @("05:49:40", "07:36", "27_EL RETIRO", 107, "LP1912"),
    @("06:15:04", "07:37", "27_EL RETIRO", 82, "LP1912"),
    @("05:49:40", "07:39", "10_OLMOS", 110, "LP1912"),
    @("07:20:40", "07:46", "16_SANTA ANA", 26, "LP1912"),
    @("06:43:40", "07:47", "14_ABASTO", 64, "LP1912"),
    @("06:15:04", "07:48", "14_ABASTO", 93, "LP1912"),
    @("06:43:40", "07:51", "215D_EL PATO", 68, "LP1912"),
    @("06:15:04", "07:52", "215D_EL PATO", 97, "LP1912"),
    @("07:47:32", "07:55", "10_OLMOS", 8, "LP1912"),
    @("07:20:40", "07:58", "16_SANTA ANA", 38, "LP1912"),
    @("07:20:40", "07:59", "23_HERNANDEZ", 39, "LP1912"),
    @("06:15:04", "08:01", "23_HERNANDEZ", 106, "LP1912"),
    @("06:43:40", "08:03", "23_HERNANDEZ", 80, "LP1912"),
    @("07:20:40", "08:03", "11_ETCHEVERRY", 43, "LP1912"),
    @("06:57:30", "08:06", "23_HERNANDEZ", 69, "LP1912"),
    @("07:47:32", "08:10", "16_SANTA ANA", 23, "LP1912"),
    @("06:15:04", "08:12", "15_ABASTO", 117, "LP1912"),
    @("07:47:32", "08:13", "10_OLMOS", 26, "LP1912"),
    @("07:47:32", "08:16", "26_HERNANDEZ", 29, "LP1912"),
    @("06:43:40", "08:21", "26_HERNANDEZ", 98, "LP1912"),
    @("06:43:40", "08:22", "16_P MOR-SANTA ANA", 99, "LP1912"),
    @("06:43:40", "08:23", "215B_EL PATO", 100, "LP1912"),
    @("06:43:40", "08:27", "84_COLONIA URQUIZA-ESC 49", 104, "LP1912"),
    @("07:47:32", "08:31", "23_HERNANDEZ", 44, "LP1912"),
    @("07:59:28", "08:33", "10_OLMOS", 34, "LP1912"),
    @("07:59:28", "08:34", "23_HERNANDEZ", 35, "LP1912"),
    @("07:59:28", "08:39", "26_HERNANDEZ", 40, "LP1912"),
    @("06:43:40", "08:42", "81_EL PELIGRO", 119, "LP1912"),
    @("07:20:40", "08:43", "14_ABASTO", 83, "LP1912"),
    @("06:57:30", "08:54", "17_ROMERO", 117, "LP1912"),
    @("08:57:13", "08:59", "11_ETCHEVERRY", 2, "LP1912"),
    @("07:20:40", "09:01", "215A_EL PATO", 101, "LP1912"),
    @("08:21:50", "09:01", "23_HERNANDEZ", 40, "LP1912"),
    @("08:57:13", "09:02", "215A_EL PATO", 5, "LP1912"),
    @("07:59:28", "09:03", "11_ETCHEVERRY", 64, "LP1912"),
    @("08:39:44", "09:04", "23_HERNANDEZ", 25, "LP1912"),
    @("08:57:13", "09:05", "23_HERNANDEZ", 8, "LP1912"),
    @("08:21:50", "09:07", "26_HERNANDEZ", 46, "LP1912"),
    @("07:20:40", "09:10", "16_P MOR-SANTA ANA", 110, "LP1912"),
    @("08:57:13", "09:11", "16_P MOR-SANTA ANA", 14, "LP1912"),
    @("08:21:50", "09:13", "10_OLMOS", 52, "LP1912"),
    @("07:20:40", "09:16", "27_EL RETIRO", 116, "LP1912"),
    @("08:57:13", "09:17", "27_EL RETIRO", 20, "LP1912"),
    @("08:21:50", "09:21", "26_HERNANDEZ", 60, "LP1912"),
    @("07:59:28", "09:22", "16_SANTA ANA", 83, "LP1912"),
    @("07:47:32", "09:22", "17_ROMERO", 95, "LP1912"),
    @("08:57:13", "09:23", "16_SANTA ANA", 26, "LP1912"),
    @("07:47:32", "09:23", "11_ETCHEVERRY", 96, "LP1912"),
    @("08:57:13", "09:24", "11_ETCHEVERRY", 27, "LP1912"),
    @("08:21:50", "09:29", "16_SANTA ANA", 68, "LP1912"),
    @("07:47:32", "09:32", "15_ABASTO", 105, "LP1912"),
    @("07:47:32", "09:33", "10_OLMOS", 106, "LP1912"),
    @("08:39:44", "09:34", "23_HERNANDEZ", 55, "LP1912"),
    @("08:39:44", "09:34", "16_SANTA ANA", 55, "LP1912"),
    @("08:57:13", "09:35", "16_SANTA ANA", 38, "LP1912"),
    @("08:57:13", "09:35", "23_HERNANDEZ", 38, "LP1912"),
    @("09:38:09", "09:41", "14_ABASTO", 3, "LP1912"),
    @("09:38:09", "09:41", "23_HERNANDEZ", 3, "LP1912"),
    @("08:21:50", "09:41", "215C_EL PATO", 80, "LP1912"),
    @("07:47:32", "09:42", "215C_EL PATO", 115, "LP1912"),
    @("07:47:32", "09:43", "14_ABASTO", 116, "LP1912"),
    @("08:57:13", "09:44", "14_ABASTO", 47, "LP1912"),
    @("09:38:09", "09:47", "16_SANTA ANA", 9, "LP1912"),
    @("08:49:51", "09:52", "15_ABASTO", 63, "LP1912"),
    @("08:49:51", "09:53", "10_OLMOS", 64, "LP1912"),
    @("09:38:09", "09:59", "16_SANTA ANA", 21, "LP1912"),
    @("09:38:09", "10:04", "11_ETCHEVERRY", 26, "LP1912"),
    @("09:38:09", "10:05", "23_HERNANDEZ", 27, "LP1912"),
    @("08:39:44", "10:06", "10_OLMOS", 87, "LP1912"),
    @("08:21:50", "10:10", "16_P MOR-SANTA ANA", 109, "LP1912"),
    @("08:57:13", "10:11", "16_P MOR-SANTA ANA", 74, "LP1912"),
    @("08:21:50", "10:12", "15_ABASTO", 111, "LP1912"),
    @("09:38:09", "10:13", "10_OLMOS", 35, "LP1912"),
    @("08:49:51", "10:20", "26_HERNANDEZ", 91, "LP1912"),
    @("08:39:44", "10:21", "26_HERNANDEZ", 102, "LP1912"),
    @("08:39:44", "10:22", "17_ROMERO", 103, "LP1912"),
    @("09:38:09", "10:24", "11_ETCHEVERRY", 46, "LP1912"),
    @("08:39:44", "10:26", "215A_EL PATO", 107, "LP1912"),
    @("08:57:13", "10:27", "215A_EL PATO", 90, "LP1912"),
    @("10:26:41", "10:33", "10_OLMOS", 7, "LP1912"),
    @("10:26:41", "10:34", "23_HERNANDEZ", 8, "LP1912"),
    @("10:26:41", "10:34", "16_SANTA ANA", 8, "LP1912"),
    @("08:49:51", "10:41", "17_ROMERO", 112, "LP1912"),
    @("08:57:13", "10:42", "17_ROMERO", 105, "LP1912"),
    @("08:49:51", "10:43", "14_ABASTO", 114, "LP1912"),
    @("08:57:13", "10:44", "14_ABASTO", 107, "LP1912"),
    @("10:26:41", "10:46", "16_SANTA ANA", 20, "LP1912"),
    @("10:26:41", "10:52", "15_ABASTO", 26, "LP1912"),
    @("10:26:41", "10:53", "10_OLMOS", 27, "LP1912"),
    @("10:26:41", "10:56", "27_EL RETIRO", 30, "LP1912"),
    @("10:56:30", "10:57", "16_SANTA ANA", 1, "LP1912"),
    @("09:38:09", "10:58", "27_EL RETIRO", 80, "LP1912"),
    @("10:26:41", "11:01", "215C_EL PATO", 35, "LP1912"),
    @("09:38:09", "11:02", "215C_EL PATO", 84, "LP1912"),
    @("10:26:41", "11:03", "11_ETCHEVERRY", 37, "LP1912"),
    @("10:26:41", "11:04", "23_HERNANDEZ", 38, "LP1912"),
    @("10:26:41", "11:06", "16_P MOR-167 Y 521", 40, "LP1912"),
    @("09:38:09", "11:07", "16_P MOR-167 Y 521", 89, "LP1912"),
    @("10:56:30", "11:08", "23_HERNANDEZ", 12, "LP1912"),
    @("10:26:41", "11:12", "15_ABASTO", 46, "LP1912"),
    @("10:26:41", "11:19", "86_EST CHICA-ESC AGRARIA", 53, "LP1912"),
    @("09:38:09", "11:20", "86_EST CHICA-ESC AGRARIA", 102, "LP1912"),
    @("11:20:07", "11:21", "16_SANTA ANA", 1, "LP1912"),
    @("09:38:09", "11:21", "26_HERNANDEZ", 103, "LP1912"),
    @("11:20:07", "11:22", "17_ROMERO", 2, "LP1912"),
    @("10:56:30", "11:24", "10_OLMOS", 28, "LP1912"),
    @("11:20:07", "11:25", "16_SANTA ANA", 5, "LP1912"),
    @("09:38:09", "11:27", "225_C ROCA-H SUR", 109, "LP1912"),
    @("09:38:09", "11:32", "81_EL PELIGRO", 114, "LP1912"),
    @("10:56:30", "11:34", "23_HERNANDEZ", 38, "LP1912"),
    @("11:20:07", "11:35", "23_HERNANDEZ", 15, "LP1912"),
    @("10:26:41", "11:35", "11_ETCHEVERRY", 69, "LP1912"),
    @("09:38:09", "11:36", "11_ETCHEVERRY", 118, "LP1912"),
    @("10:26:41", "11:41", "17_ROMERO", 75, "LP1912"),
    @("10:56:30", "11:42", "17_ROMERO", 46, "LP1912"),
    @("11:48:04", "11:49", "16_SANTA ANA", 1, "LP1912"),
    @("10:26:41", "11:51", "215B_EL PATO", 85, "LP1912"),
    @("10:56:30", "11:52", "15_ABASTO", 56, "LP1912"),
    @("11:48:04", "11:53", "16_SANTA ANA", 5, "LP1912"),
    @("10:26:41", "11:59", "225_GOMEZ", 93, "LP1912"),
    @("12:01:50", "12:01", "16_SANTA ANA", 0, "LP1912"),
    @("10:26:41", "12:02", "84_COLONIA URQUIZA-ESC 49", 96, "LP1912"),
    @("12:01:50", "12:04", "23_HERNANDEZ", 3, "LP1912"),
    @("11:20:07", "12:05", "23_HERNANDEZ", 45, "LP1912"),
    @("12:01:50", "12:06", "84_COLONIA URQUIZA-ESC 49", 5, "LP1912"),
    @("12:01:50", "12:06", "16_SANTA ANA", 5, "LP1912"),
    @("10:56:30", "12:06", "14_ABASTO", 70, "LP1912"),
    @("10:26:41", "12:06", "16_P MOR-SANTA ANA", 100, "LP1912"),
    @("11:20:07", "12:07", "16_P MOR-SANTA ANA", 47, "LP1912"),
    @("11:20:07", "12:07", "14_ABASTO", 47, "LP1912"),
    @("10:56:30", "12:10", "10_OLMOS", 74, "LP1912"),
    @("11:20:07", "12:13", "10_OLMOS", 53, "LP1912"),
    @("10:26:41", "12:14", "17_ROMERO", 108, "LP1912"),
    @("10:26:41", "12:19", "14_ABASTO", 113, "LP1912"),
    @("10:26:41", "12:20", "215A_EL PATO", 114, "LP1912"),
    @("10:56:30", "12:20", "14_ABASTO", 84, "LP1912"),
    @("11:20:07", "12:21", "14_ABASTO", 61, "LP1912"),
    @("10:26:41", "12:21", "26_HERNANDEZ", 115, "LP1912"),
    @("11:20:07", "12:21", "215A_EL PATO", 61, "LP1912"),
    @("12:01:50", "12:34", "23_HERNANDEZ", 33, "LP1912"),
    @("12:01:50", "12:34", "11_ETCHEVERRY", 33, "LP1912"),
    @("11:48:04", "12:35", "23_HERNANDEZ", 47, "LP1912"),
    @("11:48:04", "12:35", "11_ETCHEVERRY", 47, "LP1912"),
    @("10:56:30", "12:36", "27_EL RETIRO", 100, "LP1912"),
    @("11:20:07", "12:37", "27_EL RETIRO", 77, "LP1912"),
    @("12:37:14", "12:37", "23_HERNANDEZ", 0, "LP1912"),
    @("10:56:30", "12:38", "17_179 Y 38", 102, "LP1912"),
    @("10:56:30", "12:41", "10_OLMOS", 105, "LP1912"),
    @("12:37:14", "12:47", "16_SANTA ANA", 10, "LP1912"),
    @("12:01:50", "12:48", "11_ETCHEVERRY", 47, "LP1912"),
    @("11:20:07", "12:49", "11_ETCHEVERRY", 89, "LP1912"),
    @("12:37:14", "12:50", "15_ABASTO", 13, "LP1912"),
    @("12:37:14", "12:55", "10_OLMOS", 18, "LP1912"),
    @("12:01:50", "12:56", "14_ABASTO", 55, "LP1912"),
    @("11:20:07", "13:02", "15_ABASTO", 102, "LP1912"),
    @("12:37:14", "13:03", "14_ABASTO", 26, "LP1912"),
    @("12:55:01", "13:05", "23_HERNANDEZ", 10, "LP1912"),
    @("12:37:14", "13:06", "23_HERNANDEZ", 29, "LP1912"),
    @("11:20:07", "13:07", "16_P MOR-SANTA ANA", 107, "LP1912"),
    @("12:37:14", "13:08", "10_OLMOS", 31, "LP1912"),
    @("12:01:50", "13:13", "215D_EL PATO", 72, "LP1912"),
    @("11:20:07", "13:14", "215D_EL PATO", 114, "LP1912"),
    @("12:37:14", "13:15", "11_ETCHEVERRY", 38, "LP1912"),
    @("12:01:50", "13:19", "10_OLMOS", 78, "LP1912"),
    @("11:48:04", "13:20", "10_OLMOS", 92, "LP1912"),
    @("13:19:56", "13:20", "26_HERNANDEZ", 1, "LP1912"),
    @("13:19:56", "13:21", "10_OLMOS", 2, "LP1912"),
    @("11:48:04", "13:21", "26_HERNANDEZ", 93, "LP1912"),
    @("13:19:56", "13:22", "14_ABASTO", 3, "LP1912"),
    @("12:55:01", "13:23", "10_OLMOS", 28, "LP1912"),
    @("13:19:56", "13:26", "15_ABASTO", 7, "LP1912"),
    @("11:48:04", "13:27", "14_ABASTO", 99, "LP1912"),
    @("13:19:56", "13:31", "10_OLMOS", 12, "LP1912"),
    @("12:55:01", "13:32", "10_OLMOS", 37, "LP1912"),
    @("13:19:56", "13:34", "23_HERNANDEZ", 15, "LP1912"),
    @("12:55:01", "13:35", "23_HERNANDEZ", 40, "LP1912"),
    @("11:48:04", "13:36", "15_ABASTO", 108, "LP1912"),
    @("11:48:04", "13:46", "17_ROMERO", 118, "LP1912"),
    @("13:19:56", "13:46", "16_SANTA ANA", 27, "LP1912"),
    @("12:37:14", "13:47", "16_SANTA ANA", 70, "LP1912"),
    @("13:19:56", "13:50", "11_ETCHEVERRY", 31, "LP1912"),
    @("12:01:50", "13:50", "215A_EL PATO", 109, "LP1912"),
    @("12:55:01", "13:51", "11_ETCHEVERRY", 56, "LP1912"),
    @("12:37:14", "13:51", "215A_EL PATO", 74, "LP1912"),
    @("12:01:50", "13:52", "10_OLMOS", 111, "LP1912"),
    @("12:01:50", "13:55", "225_GOMEZ", 114, "LP1912"),
    @("12:37:14", "13:56", "225_GOMEZ", 79, "LP1912"),
    @("12:01:50", "13:56", "16_P MOR-167 Y 521", 115, "LP1912"),
    @("12:37:14", "13:57", "16_P MOR-167 Y 521", 80, "LP1912"),
    @("13:53:08", "13:57", "16_SANTA ANA", 4, "LP1912"),
    @("13:19:56", "14:04", "23_HERNANDEZ", 45, "LP1912"),
    @("12:37:14", "14:04", "17_ROMERO", 87, "LP1912"),
    @("13:53:08", "14:05", "23_HERNANDEZ", 12, "LP1912"),
    @("13:19:56", "14:06", "16_SANTA ANA", 47, "LP1912"),
    @("12:55:01", "14:07", "16_SANTA ANA", 72, "LP1912"),
    @("13:53:08", "14:12", "15_ABASTO", 19, "LP1912"),
    @("13:19:56", "14:16", "27_EL RETIRO", 57, "LP1912"),
    @("12:37:14", "14:17", "27_EL RETIRO", 100, "LP1912"),
    @("13:19:56", "14:19", "215C_EL PATO", 60, "LP1912"),
    @("12:37:14", "14:20", "215C_EL PATO", 103, "LP1912"),
    @("13:19:56", "14:20", "26_HERNANDEZ", 61, "LP1912"),
    @("12:37:14", "14:21", "26_HERNANDEZ", 104, "LP1912"),
    @("13:53:08", "14:28", "15_ABASTO", 35, "LP1912"),
    @("13:53:08", "14:31", "11_ETCHEVERRY", 38, "LP1912"),
    @("14:19:48", "14:34", "23_HERNANDEZ", 15, "LP1912"),
    @("14:19:48", "14:44", "15_ABASTO", 25, "LP1912"),
    @("14:19:48", "14:44", "14_ABASTO", 25, "LP1912"),
    @("13:53:08", "14:45", "14_ABASTO", 52, "LP1912"),
    @("14:19:48", "14:46", "16_SANTA ANA", 27, "LP1912"),
    @("13:19:56", "14:49", "14_ABASTO", 90, "LP1912"),
    @("12:55:01", "14:50", "14_ABASTO", 115, "LP1912"),
    @("13:19:56", "14:56", "16_P MOR-SANTA ANA", 97, "LP1912"),
    @("13:53:08", "14:57", "16_P MOR-SANTA ANA", 64, "LP1912"),
    @("13:19:56", "14:58", "215B_EL PATO", 99, "LP1912"),
    @("14:58:43", "14:58", "16_SANTA ANA", 0, "LP1912"),
    @("13:19:56", "15:00", "81_EL PELIGRO", 101, "LP1912"),
    @("14:46:52", "15:04", "23_HERNANDEZ", 18, "LP1912"),
    @("13:19:56", "15:04", "10_OLMOS", 105, "LP1912"),
    @("13:53:08", "15:05", "10_OLMOS", 72, "LP1912"),
    @("14:19:48", "15:06", "16_SANTA ANA", 47, "LP1912"),
    @("13:53:08", "15:10", "17_ROMERO", 77, "LP1912"),
    @("13:19:56", "15:13", "11_ETCHEVERRY", 114, "LP1912"),
    @("13:53:08", "15:14", "11_ETCHEVERRY", 81, "LP1912"),
    @("14:46:52", "15:20", "15_ABASTO", 34, "LP1912"),
    @("14:46:52", "15:20", "26_HERNANDEZ", 34, "LP1912"),
    @("14:19:48", "15:21", "26_HERNANDEZ", 62, "LP1912"),
    @("13:53:08", "15:28", "26_HERNANDEZ", 95, "LP1912"),
    @("14:19:48", "15:29", "14_ABASTO", 70, "LP1912"),
    @("15:31:33", "15:31", "16_SANTA ANA", 0, "LP1912"),
    @("13:53:08", "15:32", "84_COLONIA URQUIZA-ESC 49", 99, "LP1912"),
    @("14:58:43", "15:34", "23_HERNANDEZ", 36, "LP1912"),
    @("13:53:08", "15:35", "23_HERNANDEZ", 102, "LP1912"),
    @("14:46:52", "15:36", "10_OLMOS", 50, "LP1912"),
    @("13:53:08", "15:37", "10_OLMOS", 104, "LP1912"),
    @("14:19:48", "15:38", "215A_EL PATO", 79, "LP1912"),
    @("14:19:48", "15:38", "23_HERNANDEZ", 79, "LP1912"),
    @("13:53:08", "15:39", "215A_EL PATO", 106, "LP1912"),
    @("13:53:08", "15:44", "14_ABASTO", 111, "LP1912"),
    @("14:46:52", "15:44", "23_HERNANDEZ", 58, "LP1912"),
    @("14:46:52", "15:45", "14_ABASTO", 59, "LP1912"),
    @("14:19:48", "15:46", "16_P MOR-167 Y 521", 87, "LP1912"),
    @("13:53:08", "15:47", "16_P MOR-167 Y 521", 114, "LP1912"),
    @("14:19:48", "15:53", "11_ETCHEVERRY", 94, "LP1912"),
    @("15:31:33", "15:54", "16_SANTA ANA", 23, "LP1912"),
    @("14:19:48", "15:55", "17_ROMERO", 96, "LP1912"),
    @("15:31:33", "15:56", "17_ROMERO", 25, "LP1912"),
    @("14:19:48", "15:56", "27_EL RETIRO", 97, "LP1912"),
    @("15:57:48", "15:57", "17_ROMERO", 0, "LP1912"),
    @("15:57:48", "15:57", "27_EL RETIRO", 0, "LP1912"),
    @("14:46:52", "15:58", "27_EL RETIRO", 72, "LP1912"),
    @("15:31:33", "16:01", "10_OLMOS", 30, "LP1912"),
    @("15:57:48", "16:02", "16_SANTA ANA", 5, "LP1912"),
    @("15:31:33", "16:04", "23_HERNANDEZ", 33, "LP1912"),
    @("14:46:52", "16:08", "14_ABASTO", 82, "LP1912"),
    @("14:19:48", "16:14", "225_C ROCA-H SUR", 115, "LP1912"),
    @("14:58:43", "16:15", "225_C ROCA-H SUR", 77, "LP1912"),
    @("15:57:48", "16:18", "16_SANTA ANA", 21, "LP1912"),
    @("14:46:52", "16:19", "215C_EL PATO", 93, "LP1912"),
    @("14:46:52", "16:20", "26_HERNANDEZ", 94, "LP1912"),
    @("15:31:33", "16:20", "215C_EL PATO", 49, "LP1912"),
    @("14:58:43", "16:21", "26_HERNANDEZ", 83, "LP1912"),
    @("15:57:48", "16:22", "10_OLMOS", 25, "LP1912"),
    @("16:18:55", "16:26", "16_SANTA ANA", 8, "LP1912"),
    @("16:18:55", "16:28", "10_OLMOS", 10, "LP1912"),
    @("15:57:48", "16:29", "10_OLMOS", 32, "LP1912"),
    @("15:31:33", "16:29", "14_ABASTO", 58, "LP1912"),
    @("14:46:52", "16:30", "15_ABASTO", 104, "LP1912"),
    @("15:57:48", "16:34", "23_HERNANDEZ", 37, "LP1912"),
    @("15:31:33", "16:36", "11_ETCHEVERRY", 65, "LP1912"),
    @("16:18:55", "16:39", "17_ROMERO", 21, "LP1912"),
    @("15:31:33", "16:40", "17_ROMERO", 69, "LP1912"),
    @("14:46:52", "16:42", "16_P MOR-SANTA ANA", 116, "LP1912"),
    @("14:46:52", "16:42", "225_GOMEZ", 116, "LP1912"),
    @("14:58:43", "16:43", "225_GOMEZ", 105, "LP1912"),
    @("15:31:33", "16:48", "15_ABASTO", 77, "LP1912"),
    @("15:57:48", "16:50", "14_ABASTO", 53, "LP1912"),
    @("14:58:43", "16:56", "17_179 Y 38", 118, "LP1912"),
    @("16:18:55", "16:56", "10_OLMOS", 38, "LP1912"),
    @("16:18:55", "17:04", "23_HERNANDEZ", 46, "LP1912"),
    @("15:31:33", "17:04", "215A_EL PATO", 93, "LP1912"),
    @("15:57:48", "17:04", "11_ETCHEVERRY", 67, "LP1912"),
    @("16:18:55", "17:20", "26_HERNANDEZ", 62, "LP1912"),
    @("16:18:55", "17:20", "16_SANTA ANA", 62, "LP1912"),
    @("15:31:33", "17:21", "26_HERNANDEZ", 110, "LP1912"),
    @("15:31:33", "17:24", "84_COLONIA URQUIZA-ESC 49", 113, "LP1912"),
    @("16:18:55", "17:28", "14_ABASTO", 70, "LP1912"),
    @("16:18:55", "17:30", "27_EL RETIRO", 72, "LP1912"),
    @("15:57:48", "17:35", "27_EL RETIRO", 98, "LP1912"),
    @("16:18:55", "17:38", "17_ROMERO", 80, "LP1912"),
    @("16:18:55", "17:39", "215B_EL PATO", 81, "LP1912"),
    @("15:57:48", "17:40", "215B_EL PATO", 103, "LP1912"),
    @("15:57:48", "17:41", "17_ROMERO", 104, "LP1912"),
    @("15:57:48", "17:50", "16_P MOR-167 Y 521", 113, "LP1912"),
    @("15:57:48", "17:52", "81_EL PELIGRO", 115, "LP1912")
)
$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $r++
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 16:18:55"
$ws.Range("A3").Value = "Total filas: 38"
$ws.Range("A6:E43").ClearContents()
$rows = @(
    @("04:18:53", "04:45", "215A_EL PATO", 27, "LP1912"),
    @("04:40:33", "04:46", "215A_EL PATO", 6, "LP1912"),
    @("04:18:53", "05:34", "215B_EL PATO", 76, "LP1912"),
    @("04:54:03", "05:35", "215B_EL PATO", 41, "LP1912"),
    @("04:18:53", "06:11", "215A_EL PATO", 113, "LP1912"),
    @("04:54:03", "06:12", "215A_EL PATO", 78, "LP1912"),
    @("05:18:56", "06:46", "215C_EL PATO", 88, "LP1912"),
    @("04:54:03", "06:47", "215C_EL PATO", 113, "LP1912"),
    @("05:18:56", "07:11", "215A_EL PATO", 113, "LP1912"),
    @("06:15:04", "07:12", "215A_EL PATO", 57, "LP1912"),
    @("06:43:40", "07:51", "215D_EL PATO", 68, "LP1912"),
    @("06:15:04", "07:52", "215D_EL PATO", 97, "LP1912"),
    @("06:43:40", "08:23", "215B_EL PATO", 100, "LP1912"),
    @("07:20:40", "09:01", "215A_EL PATO", 101, "LP1912"),
    @("08:57:13", "09:02", "215A_EL PATO", 5, "LP1912"),
    @("08:21:50", "09:41", "215C_EL PATO", 80, "LP1912"),
    @("07:47:32", "09:42", "215C_EL PATO", 115, "LP1912"),
    @("08:39:44", "10:26", "215A_EL PATO", 107, "LP1912"),
    @("08:57:13", "10:27", "215A_EL PATO", 90, "LP1912"),
    @("10:26:41", "11:01", "215C_EL PATO", 35, "LP1912"),
    @("09:38:09", "11:02", "215C_EL PATO", 84, "LP1912"),
    @("10:26:41", "11:51", "215B_EL PATO", 85, "LP1912"),
    @("10:26:41", "12:20", "215A_EL PATO", 114, "LP1912"),
    @("11:20:07", "12:21", "215A_EL PATO", 61, "LP1912"),
    @("12:01:50", "13:13", "215D_EL PATO", 72, "LP1912"),
    @("11:20:07", "13:14", "215D_EL PATO", 114, "LP1912"),
    @("12:01:50", "13:50", "215A_EL PATO", 109, "LP1912"),
    @("12:37:14", "13:51", "215A_EL PATO", 74, "LP1912"),
    @("13:19:56", "14:19", "215C_EL PATO", 60, "LP1912"),
    @("12:37:14", "14:20", "215C_EL PATO", 103, "LP1912"),
    @("13:19:56", "14:58", "215B_EL PATO", 99, "LP1912"),
    @("14:19:48", "15:38", "215A_EL PATO", 79, "LP1912"),
    @("13:53:08", "15:39", "215A_EL PATO", 106, "LP1912"),
    @("14:46:52", "16:19", "215C_EL PATO", 93, "LP1912"),
    @("15:31:33", "16:20", "215C_EL PATO", 49, "LP1912"),
    @("15:31:33", "17:04", "215A_EL PATO", 93, "LP1912"),
    @("16:18:55", "17:39", "215B_EL PATO", 81, "LP1912"),
    @("15:57:48", "17:40", "215B_EL PATO", 103, "LP1912")
)
$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $r++
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 16:18:55"
$ws.Range("A3").Value = "Total filas: 52"
$ws.Range("A6:E57").ClearContents()
$rows = @(
    @("04:18:53", "05:43", "215A_LA PLATA", 85, "L6173"),
    @("04:40:33", "05:44", "215A_LA PLATA", 64, "L6173"),
    @("04:18:53", "06:08", "215A_LA PLATA", 110, "L6173"),
    @("04:40:33", "06:09", "215A_LA PLATA", 89, "L6173"),
    @("05:18:56", "06:32", "215C_LA PLATA", 74, "L6203"),
    @("04:40:33", "06:33", "215C_LA PLATA", 113, "L6203"),
    @("05:18:56", "06:59", "215B_LP-P MOR-1 Y 57", 101, "L6173"),
    @("06:15:04", "07:00", "215B_LP-P MOR-1 Y 57", 45, "L6173"),
    @("06:57:30", "07:01", "215B_LP-P MOR-1 Y 57", 4, "L6173"),
    @("05:49:40", "07:07", "215B_LP-P MOR-1 Y 57", 78, "L6173"),
    @("05:49:40", "07:35", "215A_LA PLATA", 106, "L6173"),
    @("07:20:40", "07:37", "215A_LA PLATA", 17, "L6173"),
    @("06:43:40", "08:06", "215C_LA PLATA", 83, "L6203"),
    @("06:15:04", "08:07", "215C_LA PLATA", 112, "L6203"),
    @("07:20:40", "08:09", "215C_LA PLATA", 49, "L6203"),
    @("06:57:30", "08:10", "215C_LA PLATA", 73, "L6203"),
    @("07:47:32", "08:13", "215C_LA PLATA", 26, "L6203"),
    @("07:59:28", "08:17", "215C_LA PLATA", 18, "L6203"),
    @("06:57:30", "08:35", "215A_LA PLATA", 98, "L6173"),
    @("07:59:28", "08:36", "215A_LA PLATA", 37, "L6173"),
    @("06:43:40", "08:38", "215A_LA PLATA", 115, "L6173"),
    @("08:21:50", "08:42", "215A_LA PLATA", 21, "L6173"),
    @("08:39:44", "08:44", "215A_LA PLATA", 5, "L6173"),
    @("07:20:40", "09:08", "215D_LA PLATA", 108, "L6203"),
    @("07:47:32", "09:09", "215D_LA PLATA", 82, "L6203"),
    @("08:21:50", "10:02", "215B_LP-P MOR-40 Y 115", 101, "L6173"),
    @("08:57:13", "10:03", "215B_LP-P MOR-40 Y 115", 66, "L6173"),
    @("08:57:13", "10:54", "215A_LA PLATA", 117, "L6173"),
    @("10:26:41", "11:13", "215C_LA PLATA", 47, "L6203"),
    @("09:38:09", "11:14", "215C_LA PLATA", 96, "L6203"),
    @("10:26:41", "12:04", "215A_LA PLATA", 98, "L6173"),
    @("12:01:50", "12:05", "215A_LA PLATA", 4, "L6173"),
    @("11:48:04", "12:06", "215A_LA PLATA", 18, "L6173"),
    @("12:01:50", "12:53", "215C_LA PLATA", 52, "L6203"),
    @("10:56:30", "12:54", "215C_LA PLATA", 118, "L6203"),
    @("12:55:01", "12:57", "215C_LA PLATA", 2, "L6203"),
    @("12:01:50", "13:30", "215B_LP-P MOR-1 Y 57", 89, "L6173"),
    @("11:48:04", "13:31", "215B_LP-P MOR-1 Y 57", 103, "L6173"),
    @("13:19:56", "14:08", "215A_LA PLATA", 49, "L6173"),
    @("12:37:14", "14:09", "215A_LA PLATA", 92, "L6173"),
    @("14:46:52", "14:47", "215D_LA PLATA", 1, "L6203"),
    @("13:19:56", "14:52", "215D_LA PLATA", 93, "L6203"),
    @("12:55:01", "14:53", "215D_LA PLATA", 118, "L6203"),
    @("13:53:08", "15:34", "215A_LA PLATA", 101, "L6173"),
    @("14:19:48", "16:13", "215C_LA PLATA", 114, "L6203"),
    @("15:31:33", "16:14", "215C_LA PLATA", 43, "L6203"),
    @("14:58:43", "16:52", "215B_LP-P MOR-40 Y 115", 114, "L6173"),
    @("15:31:33", "16:53", "215B_LP-P MOR-40 Y 115", 82, "L6173"),
    @("16:18:55", "17:13", "215A_LA PLATA", 55, "L6173"),
    @("15:31:33", "17:14", "215A_LA PLATA", 103, "L6173"),
    @("15:57:48", "17:15", "215A_LA PLATA", 78, "L6173"),
    @("16:18:55", "18:03", "215C_LA PLATA", 105, "L6203")
)
$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $r++
}
